# Updated cryptos list on Fri Feb 23 05:25:29 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "51.189.49"
Set-TextValue "E2" "  -0.71%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.957.79"
Set-TextValue "E3" "  +0.60%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "379.48"
Set-TextValue "E5" "  +0.63%  "

# Row 6 - Solana
Set-TextValue "D6" "102.40"
Set-TextValue "E6" "  -1.89%  "

# Row 7 - XRP
Set-TextValue "D7" "0.539"
Set-TextValue "E7" "  -0.63%  "

# Row 8 - USDC
Set-TextValue "E8" "  +0.06%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.590"
Set-TextValue "E9" "  +0.32%  "

# Row 10 - Avalanche
Set-TextValue "D10" "36.58"
Set-TextValue "E10" "  -1.14%  "

# Row 11 - TRON
Set-TextValue "E11" "  -0.19%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.0842"
Set-TextValue "E12" "  +0.40%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.420.01"
Set-TextValue "E13" "  +0.92%  "

# Row 14 - Chainlink
Set-TextValue "D14" "18.05"
Set-TextValue "E14" "  -2.03%  "

# Row 15 - Polkadot
Set-TextValue "D15" "7.42"
Set-TextValue "E15" "  +0.30%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.948.91"
Set-TextValue "E16" "  +0.64%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.987"
Set-TextValue "E17" "  +4.44%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "51.127.47"
Set-TextValue "E18" "  -0.72%  "

# Row 19 - ImmutableX
Set-TextValue "E19" "  -6.10%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.13"
Set-TextValue "E20" "  -2.80%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "D21" "12.52"
Set-TextValue "E21" "  -3.97%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0954"
Set-TextValue "E22" "  +0.37%  "

# Row 23 - Litecoin
Set-TextValue "D23" "68.50"
Set-TextValue "E23" "  +0.17%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "261.95"
Set-TextValue "E24" "  -0.09%  "

# Row 25 - PancakeSwap
Set-TextValue "E25" "  +1.52%  "

# Row 26 - Filecoin
Set-TextValue "D26" "8.39"
Set-TextValue "E26" "  +14.17%  "

# Row 27 - RenderToken
Set-TextValue "D27" "7.61"
Set-TextValue "E27" "  +6.56%  "

# Row 28 - Kaspa: unchanged

# Rows 29/30 swap: LEO <-> Hedera
Set-TextValue "B29" "Hedera"
Set-TextValue "C29" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D29" "0.113"
Set-TextValue "E29" "  +10.76%  "

Set-TextValue "B30" "LEO"
Set-TextValue "C30" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D30" "4.10"
Set-TextValue "E30" "  -0.66%  "

# Row 31 - Dai
Set-TextValue "E31" "  -0.04%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "25.74"
Set-TextValue "E32" "  -0.51%  "

# Row 33 - Cosmos
Set-TextValue "D33" "9.81"
Set-TextValue "E33" "  -0.27%  "

# Rows 34/35 swap: InjectiveProtocol <-> OKB
Set-TextValue "B34" "OKB"
Set-TextValue "C34" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D34" "50.51"
Set-TextValue "E34" "  -2.81%  "

Set-TextValue "B35" "InjectiveProtocol"
Set-TextValue "C35" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D35" "33.93"
Set-TextValue "E35" "  -0.69%  "

# Row 36 - Toncoin
Set-TextValue "E36" "  -2.06%  "

# Row 37 - VeChain
Set-TextValue "D37" "0.0453"
Set-TextValue "E37" "  +5.47%  "

# Row 38 - FirstDigitalUSD
Set-TextValue "E38" "  -0.03%  "

# Row 39 - LidoDAOToken
Set-TextValue "E39" "  -1.89%  "

# Row 40 - Celestia
Set-TextValue "D40" "16.87"
Set-TextValue "E40" "  -0.91%  "

# Row 41 - Stacks
Set-TextValue "D41" "2.57"

# Row 42 - Stellar
Set-TextValue "E42" "  +0.46%  "

# Row 43 - ARBITRUM
Set-TextValue "D43" "1.79"
Set-TextValue "E43" "  -2.57%  "

# Row 44 - Monero
Set-TextValue "E44" "  -2.38%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "21.24"
Set-TextValue "E45" "  -3.26%  "

# Row 46 - WEMIXToken
Set-TextValue "D46" "2.06"
Set-TextValue "E46" "  -0.19%  "

# Row 47 - TheGraph
Set-TextValue "D47" "0.275"
Set-TextValue "E47" "  +0.97%  "

# Row 48 - ApeXProtocol
Set-TextValue "E48" "  +2.17%  "

# Rows 49/50 swap: Maker <-> NEARProtocol
Set-TextValue "B49" "NEARProtocol"
Set-TextValue "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D49" "3.23"
Set-TextValue "E49" "  +1.33%  "

Set-TextValue "B50" "Maker"
Set-TextValue "C50" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D50" "2.006.45"
Set-TextValue "E50" "  -0.87%  "

# Row 51 - BEAM
Set-TextValue "D51" "0.0337"
Set-TextValue "E51" "  +3.82%  "
